$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 11; this shifts rows 11-31 down to 12-32,
# carrying all their values (and the date number format) with them.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new record.
$ws.Cells.Item(11, 1).Value = 5
$ws.Cells.Item(11, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(11, 3).Value = "Maule"
$ws.Cells.Item(11, 4).Value = 44487
$ws.Cells.Item(11, 5).Value = 7
$ws.Cells.Item(11, 6).Value = 300000000
$ws.Cells.Item(11, 7).Value = "Espárragos"
$ws.Cells.Item(11, 8).Value = "Verde"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 5000
$ws.Cells.Item(11, 11).Value = 800
$ws.Cells.Item(11, 12).Value = 800
$ws.Cells.Item(11, 13).Value = 800
$ws.Cells.Item(11, 14).Value = "`$/kilo"
$ws.Cells.Item(11, 15).Value = "Provincia de Linares"
$ws.Cells.Item(11, 16).Value = 800
$ws.Cells.Item(11, 17).Value = 1
$ws.Cells.Item(11, 18).Value = "Hortaliza"
